$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TableauBord")

# --- Row 9: add marks for BOUSSETTA Nael (A) and DEMOULIN Eddy (C), and the
#     new "Separation des controleurs" task label in column G (Lucas). The
#     date (I9) and the existing "Partage d'une liste" label (E9) are unchanged.
$ws.Range("A9").Value = "x"
$ws.Range("C9").Value = "x"
$ws.Range("G9").Value = "Separation des controleurs"

# --- Row 10: add marks for BOUSSETTA Nael (A), DEMOULIN Eddy (C) and
#     DONADONI Quentin (E), plus the new "V4" task label in column G
#     (Lucas), dated the same day as row 9.
$ws.Range("A10").Value = "x"
$ws.Range("C10").Value = "x"
$ws.Range("E10").Value = "x"
$ws.Range("G10").Value = "V4"

# Copy the date formatting from I9 onto I10 before writing the serial date
# value, so the new cell keeps the existing date number format/style instead
# of picking up a brand-new one.
$ws.Range("I9").Copy()
$ws.Range("I10").PasteSpecial(-4122)
$ws.Range("I10").Value = 44203

# Move the selection/scroll position shown when the workbook is reopened.
$ws.Activate()
$ws.Range("A8").Select()
